$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the Yoppoppo dialogue lines in column E.
#    Order matters here because the underlying shared-strings table is
#    appended to in the order cells are (re)written:
#      - E6 currently holds the lone "(Insert angry yoppoppo line here)"
#        placeholder, so rewriting it first reuses/replaces that slot.
#      - E4 and E3 currently point at the shared "Waiwai!" string (still
#        used by E2), so rewriting them creates two new shared strings.
# ---------------------------------------------------------------------------

$textE6 = @'
A Yoppoppo traveler arrives. It squints at you.
"Waiwai, Doob Belbo."
'@

$textE4 = @'
A Yoppoppo traveler arrives. It blushes and waves at you.
"Waiwai, Belbo..."
'@

$textE3 = @'
A Yoppoppo traveler arrives and waves at you.
"Waiwai, Belbo!"
'@

$ws.Range("E6").Value = $textE6
$ws.Range("E4").Value = $textE4
$ws.Range("E3").Value = $textE3

# ---------------------------------------------------------------------------
# 2) Add a new "StartingAnimation" column (F).
# ---------------------------------------------------------------------------

$ws.Range("F1").Value = "StartingAnimation"
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0

# Give the new column a sensible custom width (matches the ~18.4 char target
# as closely as this engine's column-width quantization allows).
$ws.Columns.Item(6).ColumnWidth = 17.6

# ---------------------------------------------------------------------------
# 3) Row 6 now wraps onto two lines like the other description rows, so grow
#    it to match (the other rows already carry explicit row heights).
# ---------------------------------------------------------------------------

$ws.Rows.Item(6).RowHeight = 45

# ---------------------------------------------------------------------------
# 4) Update the view state: scroll so column C is leftmost and the active
#    selection sits at F10, matching the saved workbook view.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F10").Select()
